# Backup QR Scanner data - 03/05/2025, 6:54:56 PM
# Appends a new scan-log row (row 3) to the "Scanner" sheet, mirroring the
# existing rows' layout: Number | Content | Location | Log Date | Log Time | Type

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Ftg"
$ws.Range("C3").Value = "Microbiology Lecture Hall"

# D3 holds a date string ("05/03/2025") that must stay plain text rather
# than being auto-converted into a date serial number by Excel. Force the
# cell to a text number format before assigning the value, then drop back
# to the default "Normal" style so no stray formatting is left behind.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "05/03/2025"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "6:54 PM"
$ws.Range("F3").Value = "Manual"
